$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area")

# Delete dissolved municipality "Beemster" from the options table, shifting
# the rows below it up.
$ws.Range("A57:B57").EntireRow.Delete()

# The author was working on the "area" sheet when the edit was made.
$ws.Activate()
